$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row above row 2, shifting the existing schedule down ---
$ws.Rows("2:2").Insert(-4121, 0)

# --- Drop the now-unused column G entirely (it only ever held two empty,
#     styled placeholder cells) ---
$ws.Columns("G:G").Clear()

# Re-create the visual formatting of the new row 2 by copying the row that is
# now directly below it (the old row 2, now row 3), which already carries the
# exact per-column styles we need (number/date/percent formats). This also
# seeds C2 with the correct "2025.01.12" shared-string value.
$ws.Range("A3:F3").Copy($ws.Range("A2:F2"))
$ws.Rows(2).RowHeight = 18

# --- Populate the new row 2 values ---
$ws.Range("A2").Value = 15
$ws.Range("B2").Value = 0.06805555555555555
$ws.Range("D2").Clear()
$ws.Range("E2").Value = 0.6
$ws.Range("F2").Value = "PCA"

# --- Clean up the duplicated "Actual" value that used to sit in the old D2
#     (now shifted to D3) - it does not belong there in the updated sheet ---
$ws.Range("D3").Clear()

# --- Add header "Subject" in F1, matching the style of the other headers ---
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "Subject"

# --- Update selection to match the author's final cursor position ---
$ws.Range("F2").Select()
